# PostReviewMitoWesterns.xlsx revision
# - Adds a new "SLN" measurement column (I) with a header + per-mouse values
# - Revises several existing "NDUFB8" (column C) values after re-review
# - Moves the cursor / active-cell selection to M20 (matches the saved state)
#
# (The workbook's absolute-path breadcrumb under
#  mc:AlternateContent/mc:Choice/x15ac:absPath is an Excel-managed, read-only
#  breadcrumb of the folder the file was opened from - it isn't exposed on the
#  Workbook/Application object model, so it can't be touched from COM/VBA.
#  We leave it to the host/save pipeline.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column I: "SLN" header + values -----------------------------------
$ws.Range("I1").Value = "SLN"

# --- Column C ("NDUFB8") revisions ------------------------------------------
$ws.Range("C2").Value  = 0.79560981542308729
$ws.Range("C3").Value  = 0.79440269061915325
$ws.Range("C4").Value  = 0.7774253599701495
$ws.Range("C5").Value  = 1.4333901404847211
$ws.Range("C6").Value  = 1.1991719935028888
$ws.Range("C7").Value  = 0.82551794543448231
$ws.Range("C8").Value  = 1.0338528611826157
$ws.Range("C9").Value  = 0.42272274660492465
$ws.Range("C10").Value = 0.85735468748372179
$ws.Range("C11").Value = 1.8605517592942551
$ws.Range("C12").Value = 0.87573611402298368
$ws.Range("C13").Value = 0.99277677560930144
$ws.Range("C14").Value = 1.9465817875503271
$ws.Range("C15").Value = 1.0656741128686928
$ws.Range("C16").Value = 0.80257646846846287
$ws.Range("C17").Value = 1.8101282205499289
$ws.Range("C18").Value = 1.649148909933007
$ws.Range("C19").Value = 0.80609029842317526
$ws.Range("C20").Value = 1.0356526980103515
$ws.Range("C21").Value = 2.4997768616133427

# --- Column I ("SLN") values, one per mouse row (row 7 / mouse 281 was not
#     back-filled by the author, so it intentionally has no SLN value yet) --
$ws.Range("I2").Value  = 0.43976474573583102
$ws.Range("I3").Value  = 0.4039177700424238
$ws.Range("I4").Value  = 1.1265785295055808
$ws.Range("I5").Value  = 1.4567495204272489
$ws.Range("I6").Value  = 1.5729894342889159
$ws.Range("I8").Value  = 0.57369757586764669
$ws.Range("I9").Value  = 1.0604473023786718
$ws.Range("I10").Value = 1.7238787383873213
$ws.Range("I11").Value = 0.64197638336635954
$ws.Range("I12").Value = 0.34245232816823429
$ws.Range("I13").Value = 0.5512732208312906
$ws.Range("I14").Value = 2.7937192411456748
$ws.Range("I15").Value = 2.5572209023274066
$ws.Range("I16").Value = 1.2215606663616907
$ws.Range("I17").Value = 0.60611960391135367
$ws.Range("I18").Value = 1.1927796215804063
$ws.Range("I19").Value = 1.422975280265582
$ws.Range("I20").Value = 1.1447176216445916
$ws.Range("I21").Value = 0.60384917063911503

# --- Match the saved selection/active cell ----------------------------------
$ws.Range("M20").Select()
